# Update column G ("K") values for rows 2-34 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(5, 1, 3, 8, 6, 3, 4, 8, 1, 4, 7, 1, 0, 2, 9, 5, 5, 12, 9, 2, 8, 2, 4, 8, 4, 4, 7, 8, 6, 4, 2, 3, 2)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 7).Value = $v
    $row++
}
